$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column stores plain text (e.g. "62.881.32", "1.00",
# "0.0000170") rather than numbers, so pre-format the cells that are
# about to receive new values as Text. This stops Excel from
# auto-coercing numeric-looking strings and losing formatting such as
# trailing/leading zeros. (Each contiguous block is formatted with its
# own statement - this engine only honors the first area of a
# multi-area "A1,A2,.." Range() when setting NumberFormat.)
$ws.Range("D2:D7").NumberFormat = "@"
$ws.Range("D9:D10").NumberFormat = "@"
$ws.Range("D13:D14").NumberFormat = "@"
$ws.Range("D16:D24").NumberFormat = "@"
$ws.Range("D26:D29").NumberFormat = "@"
$ws.Range("D31:D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36:D45").NumberFormat = "@"
$ws.Range("D47:D51").NumberFormat = "@"

$ws.Range("D2").Value = "62.881.32"
$ws.Range("E2").Value = "  +2.01%  "

$ws.Range("D3").Value = "3.469.46"
$ws.Range("E3").Value = "  +2.14%  "

$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").Value = "577.22"
$ws.Range("E5").Value = "  +0.08%  "

$ws.Range("D6").Value = "147.73"
$ws.Range("E6").Value = "  +3.53%  "

$ws.Range("D7").Value = "3.472.53"
$ws.Range("E7").Value = "  +2.22%  "

$ws.Range("E8").Value = "  -0.02%  "

$ws.Range("D9").Value = "0.479"
$ws.Range("E9").Value = "  +1.31%  "

$ws.Range("D10").Value = "7.65"
$ws.Range("E10").Value = "  +0.05%  "

$ws.Range("E11").Value = "  +1.29%  "

$ws.Range("E12").Value = "  +4.34%  "

$ws.Range("D13").Value = "4.065.29"
$ws.Range("E13").Value = "  +2.27%  "

$ws.Range("D14").Value = "29.78"
$ws.Range("E14").Value = "  +6.46%  "

$ws.Range("E15").Value = "  +2.81%  "

$ws.Range("D16").Value = "3.474.86"
$ws.Range("E16").Value = "  +2.62%  "

$ws.Range("D17").Value = "0.0000170"
$ws.Range("E17").Value = "  +0.07%  "

$ws.Range("D18").Value = "62.920.47"
$ws.Range("E18").Value = "  +2.00%  "

$ws.Range("D19").Value = "6.33"
$ws.Range("E19").Value = "  +3.50%  "

$ws.Range("D20").Value = "14.37"
$ws.Range("E20").Value = "  +5.42%  "

$ws.Range("D21").Value = "9.21"
$ws.Range("E21").Value = "  +1.03%  "

$ws.Range("D22").Value = "389.10"
$ws.Range("E22").Value = "  +0.13%  "

$ws.Range("D23").Value = "0.557"
$ws.Range("E23").Value = "  +1.68%  "

$ws.Range("D24").Value = "74.65"
$ws.Range("E24").Value = "  -0.03%  "

$ws.Range("E25").Value = "  +0.07%  "

$ws.Range("D26").Value = "3.607.64"
$ws.Range("E26").Value = "  +2.10%  "

$ws.Range("D27").Value = "0.0000114"
$ws.Range("E27").Value = "  +1.04%  "

$ws.Range("D28").Value = "0.179"
$ws.Range("E28").Value = "  -1.09%  "

$ws.Range("D29").Value = "7.53"
$ws.Range("E29").Value = "  +2.02%  "

$ws.Range("E30").Value = "  +0.26%  "

$ws.Range("D31").Value = "8.13"
$ws.Range("E31").Value = "  +1.95%  "

$ws.Range("D32").Value = "2.12"
$ws.Range("E32").Value = "  -0.90%  "

$ws.Range("E33").Value = "  +0.04%  "

$ws.Range("D34").Value = "1.35"
$ws.Range("E34").Value = "  -3.39%  "

$ws.Range("E35").Value = "  +1.23%  "

$ws.Range("B36").Value = "EnergySwap"
$ws.Range("C36").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D36").Value = "32.15"
$ws.Range("E36").Value = "  +20.25%  "

$ws.Range("B37").Value = "NEARProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D37").Value = "5.25"
$ws.Range("E37").Value = "  +2.87%  "

$ws.Range("B38").Value = "Aptos"
$ws.Range("C38").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D38").Value = "7.04"
$ws.Range("E38").Value = "  +1.44%  "

$ws.Range("D39").Value = "170.49"
$ws.Range("E39").Value = "  +1.25%  "

$ws.Range("D40").Value = "1.56"
$ws.Range("E40").Value = "  +5.65%  "

$ws.Range("D41").Value = "3.508.98"
$ws.Range("E41").Value = "  +2.33%  "

$ws.Range("D42").Value = "0.0753"
$ws.Range("E42").Value = "  -1.29%  "

$ws.Range("D43").Value = "0.799"
$ws.Range("E43").Value = "  +2.08%  "

$ws.Range("D44").Value = "42.37"
$ws.Range("E44").Value = "  -0.10%  "

$ws.Range("D45").Value = "4.45"
$ws.Range("E45").Value = "  +0.33%  "

$ws.Range("E46").Value = "  +2.44%  "

$ws.Range("D47").Value = "1.20"
$ws.Range("E47").Value = "  +3.74%  "

$ws.Range("D48").Value = "2.612.72"
$ws.Range("E48").Value = "  +5.58%  "

$ws.Range("D49").Value = "2.26"
$ws.Range("E49").Value = "  +11.02%  "

$ws.Range("D50").Value = "22.90"
$ws.Range("E50").Value = "  +1.04%  "

$ws.Range("D51").Value = "6.71"
$ws.Range("E51").Value = "  +0.90%  "
